$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "SinnaSone" -> "Sinna Sone" (row 10, column C)
$ws.Range("C10").Value = "Sinna Sone"

# Update Dung's task: was "Xây dựng tài liệu hướng dẫn cài đặt, vận hành",
# now "Xây dựng README.md, LICENSE" (row 11, column B)
$ws.Range("B11").Value = "Xây dựng README.md, LICENSE"

# The new text is shorter and no longer needs two wrapped lines, so the
# row shrinks from its previous auto-fit height back down to a single line.
$ws.Rows.Item(11).RowHeight = 18

# Reflect where the user's selection ended up after the edit.
$ws.Range("B5").Select()
